$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing covid_deaths (column C) values ---
$updates = @{
    595  = 5
    904  = 21
    920  = 25
    925  = 30
    930  = 18
    945  = 26
    969  = 9
    974  = 37
    978  = 18
    986  = 36
    993  = 4
    1020 = 42
    1025 = 49
    1030 = 39
    1040 = 45
    1049 = 20
    1050 = 37
    1053 = 12
    1055 = 51
    1057 = 18
    1059 = 46
    1061 = 6
    1062 = 16
    1063 = 25
    1064 = 35
    1070 = 43
    1073 = 20
    1074 = 40
    1079 = 17
    1080 = 37
    1083 = 8
    1085 = 46
    1090 = 33
    1091 = 2
    1094 = 8
    1095 = 12
    1096 = 37
    1098 = 5
    1099 = 9
    1100 = 13
    1101 = 43
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# --- Row 1097 age group changed from "0-19" to "40-49" ---
$ws.Cells.Item(1097, 2).Value = "40-49"

# --- Append new rows 1102-1120 ---
$newRows = @(
    @(44183, "50-59", 3),
    @(44183, "60-69", 9),
    @(44183, "70-79", 21),
    @(44183, "80+", 33),
    @(44184, "50-59", 3),
    @(44184, "60-69", 7),
    @(44184, "70-79", 17),
    @(44184, "80+", 33),
    @(44185, "50-59", 2),
    @(44185, "60-69", 7),
    @(44185, "70-79", 17),
    @(44185, "80+", 29),
    @(44186, "50-59", 2),
    @(44186, "60-69", 9),
    @(44186, "70-79", 19),
    @(44186, "80+", 27),
    @(44187, "60-69", 2),
    @(44187, "70-79", 4),
    @(44187, "80+", 4)
)

$r = 1102
foreach ($entry in $newRows) {
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $r = $r + 1
}
